$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1. Structural cleanup: remove the four horizontal-rule paragraphs
#    (decorative <w:pict> rectangles used as <hr> separators).
#    Delete from bottom to top so earlier paragraph indices stay valid.
# --------------------------------------------------------------------------
$hrIndexes = @(21, 17, 6, 2)
foreach ($idx in $hrIndexes) {
    $p = $d.Paragraphs($idx)
    $p.Range.Delete()
}

# --------------------------------------------------------------------------
# 2. The paragraph that used to follow the final horizontal rule ("This post
#    is worth ... completion credit. Due Sunday by 11:59pm.") now starts the
#    flow directly, so its style changes from "First Paragraph" to
#    "Body Text". Do this before re-applying run-level formatting, because
#    assigning .Style resets direct run formatting on that paragraph.
# --------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Style = "Body Text"

# --------------------------------------------------------------------------
# 3. Point value correction: 5 -> 10
# --------------------------------------------------------------------------
$d.Content.Find.ClearFormatting()
$replaced = $d.Content.Find.Execute("This post is worth 5 points completion credit.", $true, $false, $false, $false, $false, $true, 1, $false, "This post is worth 10 points completion credit.", 2)

# --------------------------------------------------------------------------
# 4. Apply the "Inter" font across the whole document body.
# --------------------------------------------------------------------------
$d.Content.Font.Name = "Inter"

# --------------------------------------------------------------------------
# 5. Color the bold "label" runs (headings/prompts/callouts) with the
#    accent color #0F4761. Each lookup is scoped to the paragraph that
#    contains it so short/ambiguous labels (e.g. "two") match only the
#    intended run.
# --------------------------------------------------------------------------
function Set-LabelColor($paraIndex, $label) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range.Duplicate
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Color = 6375183
        $rng.Font.Bold = $true
        $rng.Font.BoldBi = $true
    }
}

Set-LabelColor 1  "NotebookLM Discussion Post - Week 1"
Set-LabelColor 3  "Boaler, Chapter 3: A Mindset for Mathematics"
Set-LabelColor 4  "Ambitious Science Teaching, Chapter 1: How Do You Begin the Year?"
Set-LabelColor 8  "Prompt 1: What NotebookLM tool did you try?"
Set-LabelColor 10 "Prompt 2: How did you use it?"
Set-LabelColor 12 "Prompt 3: What’s one idea from the reading that will change how you teach?"
Set-LabelColor 16 "two"
Set-LabelColor 16 "2-3 sentences"
Set-LabelColor 17 "NOT accepted:"
Set-LabelColor 18 "This post is worth 10 points completion credit."

# --------------------------------------------------------------------------
# 6. Tighten page margins to 0.5" (720 twips / 36pt) on all sides.
# --------------------------------------------------------------------------
$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36
$d.PageSetup.RightMargin = 36

Write-Host "Edit complete."
